$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 100
$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 300

$ws.Range("B5").Select()
